$d = $word.ActiveDocument

# 1) Update the title text at both occurrences.
$d.Content.Find.Execute(
    "Grammar input set model specificatíon (layers kinds).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Grammar input set model specificatíon (Statement layer kinds).", 2
) | Out-Null

$newTitle = "Grammar input set model specificatíon (Statement layer kinds)."
$newBodyText = "Dimensional input set model specificatíon (from Statement layer, ordered SPOs: order criteria, comparisons. Kinds / class / occurrence / instance order criteria?). Value, Previous, Distance, Next. Dimension, Unit, Measure, Value (aggregated ordered statements)."

# 2) Find every paragraph that now holds the updated title text.
$count = $d.Paragraphs.Count
$targets = @()
for ($i = 1; $i -le $count; $i++) {
    $paraText = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($paraText -eq $newTitle) {
        $targets += $i
    }
}

# 3) Walk the hits back-to-front so earlier indices stay valid while we insert.
[array]::Reverse($targets)

ForEach ($idx in $targets) {
    $titlePara = $d.Paragraphs($idx)

    # The existing paragraph right after the title is the (already present)
    # blank spacer paragraph - insert the two new paragraphs right in front
    # of it, preserving its own position/content.
    $spacerPara = $titlePara.Next()

    $blankRange = $spacerPara.Range
    $blankRange.Collapse(1)
    $blankRange.InsertParagraphBefore() | Out-Null

    $spacerPara2 = $titlePara.Next().Next()
    $bodyRange = $spacerPara2.Range
    $bodyRange.Collapse(1)
    $bodyRange.InsertParagraphBefore() | Out-Null

    $newBodyPara = $titlePara.Next().Next()
    $newBodyRange = $newBodyPara.Range
    $newBodyRange.Collapse(1)
    $newBodyRange.InsertAfter($newBodyText)

    $newBodyRange = $newBodyPara.Range
    $newBodyRange.Font.Size = 8
    $newBodyRange.Font.SizeBi = 8
}
